# Fill in the previously-empty "average / sum" rows for the S1 block
# (rows 29-30), mirroring the existing S2 block rows (60-61), and rename
# the "Ecart moyen" / "Somme des écarts" labels (FR) to their Dutch
# equivalents "Gemiddeld verschil" / "Som der verschillen" everywhere
# they are used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/number formats/borders) of the existing
# S2 summary rows onto the empty gap at rows 29-30 so the new cells pick
# up identical styles (s="5" / s="8" / s="13" in the saved file).
$ws.Range("A60:B61").Copy()
$ws.Range("A29:B30").PasteSpecial(-4122)

# New "Gemiddeld verschil" (average) row for the first (S1) block.
$ws.Cells.Item(29, 1).Value2 = "Gemiddeld verschil"
$ws.Cells.Item(29, 2).Formula = "=SUM(B2:B28)/22"

# New "Som der verschillen" (sum) row for the first (S1) block.
$ws.Cells.Item(30, 1).Value2 = "Som der verschillen"
$ws.Cells.Item(30, 2).Formula = "=SUM(B2:B28)"

# Relabel the existing S2-block summary rows (was French, now Dutch) so
# both blocks reuse the very same shared strings and the old French
# strings fall out of the shared-string table entirely.
$ws.Cells.Item(60, 1).Value2 = "Gemiddeld verschil"
$ws.Cells.Item(61, 1).Value2 = "Som der verschillen"

# Widen the second chart ("Graphique 2", S2 trend chart) so its right
# edge extends further to the right (column H -> column L).
$co2 = $ws.ChartObjects().Item(2)
$co2.Left = 396.072265625
$co2.Width = 550.8749
